# Apply updates to "Översikt VADSTENA" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-16: 46070 -> 46072
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46072
}

# Row 11's record (A 58926-2025) moved up to row 8, shifting rows 8-10 down to 9-11.
# Capture the data that needs to move (old row 8, 9, 10, 11) before overwriting.
$oldA8 = $ws.Cells.Item(8, 1).Value2
$oldB8 = $ws.Cells.Item(8, 2).Value2
$oldG8 = $ws.Cells.Item(8, 7).Value2

$oldA9 = $ws.Cells.Item(9, 1).Value2
$oldB9 = $ws.Cells.Item(9, 2).Value2
$oldG9 = $ws.Cells.Item(9, 7).Value2

$oldA10 = $ws.Cells.Item(10, 1).Value2
$oldB10 = $ws.Cells.Item(10, 2).Value2
$oldG10 = $ws.Cells.Item(10, 7).Value2

$oldA11 = $ws.Cells.Item(11, 1).Value2
$oldB11 = $ws.Cells.Item(11, 2).Value2
$oldG11 = $ws.Cells.Item(11, 7).Value2

# New row 8 = old row 11
$ws.Cells.Item(8, 1).Value = $oldA11
$ws.Cells.Item(8, 2).Value = $oldB11
$ws.Cells.Item(8, 7).Value = $oldG11

# New row 9 = old row 8
$ws.Cells.Item(9, 1).Value = $oldA8
$ws.Cells.Item(9, 2).Value = $oldB8
$ws.Cells.Item(9, 7).Value = $oldG8

# New row 10 = old row 9
$ws.Cells.Item(10, 1).Value = $oldA9
$ws.Cells.Item(10, 2).Value = $oldB9
$ws.Cells.Item(10, 7).Value = $oldG9

# New row 11 = old row 10
$ws.Cells.Item(11, 1).Value = $oldA10
$ws.Cells.Item(11, 2).Value = $oldB10
$ws.Cells.Item(11, 7).Value = $oldG10
